$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 191; all existing rows from 191 downward
# shift down by 2 (191->193 ... 269->271), carrying their formatting along.
$ws.Rows.Item(191).Resize(2).Insert()

# Populate new row 191
$ws.Cells.Item(191, 1).Value = 9
$ws.Cells.Item(191, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(191, 3).Value = "Metropolitana"
$ws.Cells.Item(191, 4).Value = 44489
$ws.Cells.Item(191, 5).Value = 13
$ws.Cells.Item(191, 6).Value = 100112031
$ws.Cells.Item(191, 7).Value = "Poroto verde"
$ws.Cells.Item(191, 8).Value = "Magnum"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 34
$ws.Cells.Item(191, 11).Value = 39000
$ws.Cells.Item(191, 12).Value = 40000
$ws.Cells.Item(191, 13).Value = 39500
$ws.Cells.Item(191, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(191, 15).Value = "Perú"
$ws.Cells.Item(191, 16).Value = 1580
$ws.Cells.Item(191, 17).Value = 25
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# Populate new row 192
$ws.Cells.Item(192, 1).Value = 9
$ws.Cells.Item(192, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(192, 3).Value = "Metropolitana"
$ws.Cells.Item(192, 4).Value = 44489
$ws.Cells.Item(192, 5).Value = 13
$ws.Cells.Item(192, 6).Value = 100112031
$ws.Cells.Item(192, 7).Value = "Poroto verde"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 16
$ws.Cells.Item(192, 11).Value = 40000
$ws.Cells.Item(192, 12).Value = 43000
$ws.Cells.Item(192, 13).Value = 41500
$ws.Cells.Item(192, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(192, 15).Value = "Perú"
$ws.Cells.Item(192, 16).Value = 1660
$ws.Cells.Item(192, 17).Value = 25
$ws.Cells.Item(192, 18).Value = "Hortaliza"
